$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (shared strings used by row 1 headers) ---
$ws.Range("F1").Value = "Az_pecentage_not_matched"
$ws.Range("H1").Value = "Asctb_percentage_not_matched"

# --- Update data values (rows 2-7) ---
# Row 2: lung
$ws.Range("C2").Value = 67
$ws.Range("D2").Value = 102
$ws.Range("F2").Value = 67.164179104477611
$ws.Range("H2").Value = 23.52941176470588

# Row 3: pancreas
$ws.Range("C3").Value = 13
$ws.Range("D3").Value = 32
$ws.Range("F3").Value = 23.07692307692308
$ws.Range("H3").Value = 37.5

# Row 4: kidney
$ws.Range("C4").Value = 63
$ws.Range("D4").Value = 63
$ws.Range("F4").Value = 26.984126984126981
$ws.Range("H4").Value = 17.460317460317459

# Row 5: brain
$ws.Range("C5").Value = 190
$ws.Range("D5").Value = 127
$ws.Range("F5").Value = 99.473684210526315

# Row 6: bone_marrow
$ws.Range("C6").Value = 49
$ws.Range("D6").Value = 55
$ws.Range("F6").Value = 18.367346938775508
$ws.Range("H6").Value = 23.63636363636364

# Row 7: blood_pmbc
$ws.Range("C7").Value = 78
$ws.Range("D7").Value = 32
$ws.Range("F7").Value = 15.38461538461539
$ws.Range("H7").Value = 43.75
